$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individual Costs")

$ws.Range("C2").Value = 242.9301789682836
$ws.Range("C3").Value = 241.7533522765146
$ws.Range("C4").Value = 244.4072642249967
$ws.Range("C5").Value = 230.6829339348994
$ws.Range("C6").Value = 219.1516994869718
$ws.Range("C7").Value = 211.218754640267
$ws.Range("C8").Value = 203.3390079515765
$ws.Range("C9").Value = 201.3797721673737
$ws.Range("C10").Value = 245.4722048135623
$ws.Range("C11").Value = 254.8313851478964
$ws.Range("C12").Value = 260.2003103654174
$ws.Range("C13").Value = 259.9537246242008
$ws.Range("C14").Value = 252.9227432672594
$ws.Range("C15").Value = 267.3152168258129
$ws.Range("C16").Value = 300.3356876444313
$ws.Range("C17").Value = 305.7807857621891
$ws.Range("C18").Value = 258.7384440859091
$ws.Range("C19").Value = 237.7027714931847
$ws.Range("C20").Value = 281.3168672212262
$ws.Range("C21").Value = 256.801896043078
$ws.Range("C22").Value = 238.1904045372331
$ws.Range("C23").Value = 229.3050071278335
$ws.Range("C24").Value = 203.1075069825606
$ws.Range("C25").Value = 217.5450755808948
